# Quarterly indexing esoteric bug-fix operation
#
# Each data row (rows 2-16) held a series of q-o-q error values starting
# in column B. The indexing was off by one quarter: a new leading value
# needs to be inserted at the front of every row (column B) and the rest
# of that row's existing values shift one column to the right, with the
# series staying capped at column K (i.e. the trailing value that would
# overflow past K is simply dropped / not carried forward).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full row values (after the fix), keyed by row number, left to right
# starting at column B.
$newRowValues = @{
    2  = @(-0.6603092772102132, -0.1340702663479956, -0.2758992130872696, 0.3670754083147943, 0.4914881146443768, 0.09385790860675949, 0.01768919490591373, 0.2940112530288354, 0.1227914307171113, 0.2759388362258526)
    3  = @(-0.15162438770796, -0.293453334447234, 0.3495212869548299, 0.4739339932844123, 0.07630378724679503, 0.0001350735459492769, 0.2764571316688709, 0.1052373093571469, 0.2583847148658881, -0.01904405955723064)
    4  = @(-0.2053460154962278, 0.4376286059058361, 0.5620413122354185, 0.1644111061978012, 0.08824239249695551, 0.3645644506198771, 0.1933446283081531, 0.3464920338168943, 0.06906325939377558, -0.09392443396517081)
    5  = @(0.6162032393936197, 0.7406159457232021, 0.3429857396855849, 0.2668170259847391, 0.5431390841076607, 0.3719192617959367, 0.525066667304678, 0.2476378928815592, 0.0846501995226128, -0.03943237587190501)
    6  = @(1.652643173475852, 1.255012967438235, 1.178844253737389, 1.455166311860311, 1.283946489548587, 1.437093895057328, 1.159665120634209, 0.9966774272752628, 0.8725948518807449, 1.388747888886706)
    7  = @(0.3110387314724781, 0.2348700177716323, 0.5111920758945538, 0.3399722535828299, 0.4931196590915711, 0.2156908846684524, 0.05270319130950599, -0.07137938408501182, 0.444773652920949, 0.2348700177716323)
    8  = @(0.2388379152847414, 0.5151599734076631, 0.343940151095939, 0.4970875566046802, 0.2196587821815615, 0.0566710888226151, -0.06741148657190271, 0.4487415504340581, 0.2388379152847414)
    9  = @(0.6508000635779043, 0.4795802412661804, 0.6327276467749217, 0.3552988723518029, 0.1923111789928565, 0.06822860359833866, 0.5843816406042994, 0.3744780054549828)
    10 = @(0.2387740594105157, 0.3919214649192569, 0.1144926904961382, -0.04849500286280822, -0.172577578257326, 0.3435754587486348, 0.1336718235993181)
    11 = @(0.3465902496671606, 0.0691614752440418, -0.09382621811490459, -0.2179087935094224, 0.2982442434965384, 0.08834060834722172)
    12 = @(0.00230005330798793, -0.1606876400509585, -0.2847702154454763, 0.2313828215604846, 0.02147918641116785)
    13 = @(-0.1902738424076751, -0.3143564178021929, 0.201796619203768, -0.00810701594554874)
    14 = @(-0.3325070745318338, 0.1836459624741271, -0.02625767267518964)
    15 = @(0.1656141382254278, -0.04428949692388896)
    16 = @(-0.09587373626955231)
}

foreach ($rowNum in $newRowValues.Keys) {
    $values = $newRowValues[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $colNum = 2 + $i   # column B = 2
        $ws.Cells.Item($rowNum, $colNum).Value = $values[$i]
    }
}
